$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 119, shifting rows 119:149 down to 120:150
$ws.Rows.Item(119).Insert()

# Populate the new row 119 with this week's data point
$ws.Range("A119").Value = 4
$ws.Range("B119").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C119").Value = "Los Lagos"
$ws.Range("D119").Value = 45135
$ws.Range("E119").Value = 10
$ws.Range("F119").Value = 100112026
$ws.Range("G119").Value = "Haba"
$ws.Range("H119").Value = "Sin especificar"
$ws.Range("I119").Value = "Primera"
$ws.Range("J119").Value = 80
$ws.Range("K119").Value = 18000
$ws.Range("L119").Value = 18000
$ws.Range("M119").Value = 18000
$ws.Range("N119").Value = "$/saco 25 kilos"
$ws.Range("O119").Value = "Provincia de Limarí"
$ws.Range("P119").Value = 720
$ws.Range("Q119").Value = 25
$ws.Range("R119").Value = "Hortaliza"
